$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (K = strikeouts) values per regenerated save_data
$kValues = @{
    2 = 2
    3 = 0
    4 = 2
    5 = 3
    6 = 0
    7 = 1
    8 = 0
    9 = 2
    10 = 3
    11 = 1
    12 = 0
    13 = 2
    14 = 2
    15 = 1
    16 = 1
    17 = 0
    18 = 3
    19 = 4
    20 = 3
    21 = 3
    22 = 1
    23 = 1
    24 = 2
    25 = 2
    26 = 0
    27 = 0
    28 = 1
    29 = 0
    30 = 1
    31 = 1
    32 = 2
    33 = 3
    34 = 1
    35 = 1
    36 = 0
    37 = 2
    38 = 2
    39 = 2
    40 = 2
    41 = 1
    42 = 1
    43 = 3
    44 = 0
    45 = 3
    46 = 0
    47 = 2
    48 = 1
    49 = 1
    50 = 2
    51 = 3
    52 = 1
    53 = 1
    54 = 1
    55 = 0
    56 = 1
    57 = 2
    58 = 0
    59 = 1
    60 = 2
    61 = 1
    62 = 1
    63 = 0
    64 = 2
    65 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Host "Updated $($kValues.Count) K values in column G"
